$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 85322.09149999999
$ws.Range("B2").Value = 65000
$ws.Range("G2").Value = 20000
$ws.Range("H2").Value = 14782.18299999999
$ws.Range("A3").Value = 66643.49400000001
$ws.Range("B3").Value = 55000
$ws.Range("G3").Value = 10000
$ws.Range("H3").Value = 6299.988000000005
$ws.Range("A4").Value = 63850.3595
$ws.Range("B4").Value = 45000
$ws.Range("C4").Value = 25000
$ws.Range("D4").Value = 20000
$ws.Range("E4").Value = 10000
$ws.Range("F4").Value = 10000
$ws.Range("H4").Value = 13599.719
$ws.Range("A5").Value = 61979.8705
$ws.Range("B5").Value = 45000
$ws.Range("C5").Value = 25000
$ws.Range("D5").Value = 20000
$ws.Range("E5").Value = 10000
$ws.Range("F5").Value = 10000
$ws.Range("H5").Value = 11785.74099999999
$ws.Range("A6").Value = 62565.603
$ws.Range("H6").Value = 24922.20600000001
$ws.Range("A7").Value = 74491.1695
$ws.Range("H7").Value = 46802.33900000001
$ws.Range("A8").Value = 60588.8635
$ws.Range("H8").Value = 32460.727
$ws.Range("A9").Value = 67205.92599999999
$ws.Range("H9").Value = 38156.85199999999
$ws.Range("A10").Value = 75170.359
$ws.Range("B10").Value = 22500
$ws.Range("C10").Value = 12500
$ws.Range("D10").Value = 10000
$ws.Range("E10").Value = 5000
$ws.Range("F10").Value = 5000
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 44860.71799999999
$ws.Range("A11").Value = 82295.1295
$ws.Range("B11").Value = 32500
$ws.Range("C11").Value = 12500
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 36817.25899999999
$ws.Range("A12").Value = 83060.50599999999
$ws.Range("B12").Value = 55000
$ws.Range("G12").Value = 10000
$ws.Range("H12").Value = 13181.01199999999
$ws.Range("A13").Value = 79680.716
$ws.Range("H13").Value = -653.5679999999993
$ws.Range("A14").Value = 87679.5435
$ws.Range("H14").Value = 7141.087
$ws.Range("A15").Value = 85245.04149999999
$ws.Range("H15").Value = 4771.082999999984
$ws.Range("A16").Value = 84102.6715
$ws.Range("H16").Value = 3643.342999999993
$ws.Range("A17").Value = 82366.4325
$ws.Range("H17").Value = 1632.864999999991
$ws.Range("A18").Value = 78912.772
$ws.Range("H18").Value = -2086.456000000006
$ws.Range("A19").Value = 77523.92999999999
$ws.Range("H19").Value = -3163.140000000007
$ws.Range("A20").Value = 79823.067
$ws.Range("H20").Value = -495.8660000000091
$ws.Range("A21").Value = 77799.3585
$ws.Range("H21").Value = -635.2829999999958
$ws.Range("A22").Value = 77147.185
$ws.Range("H22").Value = 537.3699999999953
$ws.Range("A23").Value = 75950.549
$ws.Range("H23").Value = 1847.097999999998
$ws.Range("A24").Value = 73885.5585
$ws.Range("H24").Value = 2595.116999999998
$ws.Range("A25").Value = 71567.1335
$ws.Range("H25").Value = 1349.266999999993
